$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Inventory)
$ws.Range("B4").Value = 231000000.0
$ws.Range("C4").Value = 192000000.0
$ws.Range("D4").Value = 160000000.0
$ws.Range("E4").Value = 166000000.0
$ws.Range("F4").Value = 148000000.0

# Row 14 (Accounts Payable)
$ws.Range("B14").Value = 25000000.0
$ws.Range("C14").Value = 30000000.0
$ws.Range("D14").Value = 37000000.0
$ws.Range("E14").Value = 52000000.0
$ws.Range("F14").Value = 27000000.0

# Row 22 (Long Term Tax Liability (Deferred))
$ws.Range("B22").Value = -523000000.0
$ws.Range("C22").Value = -498000000.0
$ws.Range("D22").Value = -464000000.0
$ws.Range("E22").Value = -430000000.0
$ws.Range("F22").Value = -407000000.0
